$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D: "digikey"
$ws.Range("D1").Value = "digikey"

# New expense value for Carlos Mariscal (row 2) under digikey
$ws.Range("D2").Value = 33.83

# Reflect the active selection recorded in the saved file
$ws.Range("D2").Select()
